$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric but must stay as text
# (matching the source data which stores them as strings, e.g. "1.000").
# Temporarily force Text format so Excel does not auto-convert them to numbers,
# then restore the original "Normal" style so no visual/style change remains.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.158.22"
$ws.Range("E2").Value = "  +5.95%  "

$ws.Range("D3").Value = "1.916.06"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  -0.77%  "

$ws.Range("D5").Value = "329.81"
$ws.Range("E5").Value = "  +4.64%  "

$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.78%  "

$ws.Range("E7").Value = "  +3.17%  "

$ws.Range("D8").Value = "0.4076"
$ws.Range("E8").Value = "  +4.32%  "

$ws.Range("D9").Value = "0.08520"
$ws.Range("E9").Value = "  +2.40%  "

$ws.Range("D10").Value = "42.89"
$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("E11").Value = "  +1.84%  "

$ws.Range("D12").Value = "22.41"
$ws.Range("E12").Value = "  +10.48%  "

$ws.Range("E13").Value = "  +4.03%  "

$ws.Range("D14").Value = "1.913.84"
$ws.Range("E14").Value = "  +2.49%  "

$ws.Range("D15").Value = "7.402"
$ws.Range("E15").Value = "  +2.35%  "

$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("E17").Value = "  +4.23%  "

$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").Value = "0.06697"
$ws.Range("E19").Value = "  -0.45%  "

$ws.Range("D20").Value = "18.38"
$ws.Range("E20").Value = "  +4.28%  "

$ws.Range("D21").Value = "0.9995"
$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("D22").Value = "6.005"
$ws.Range("E22").Value = "  +1.72%  "

$ws.Range("D23").Value = "30.166.29"
$ws.Range("E23").Value = "  +5.78%  "

$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  +2.26%  "

$ws.Range("D25").Value = "2.210"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "2.134.99"
$ws.Range("E26").Value = "  +2.70%  "

$ws.Range("D27").Value = "161.06"
$ws.Range("E27").Value = "  +2.14%  "

$ws.Range("D28").Value = "21.11"
$ws.Range("E28").Value = "  +2.81%  "

$ws.Range("E29").Value = "  +0.25%  "

$ws.Range("D30").Value = "128.87"
$ws.Range("E30").Value = "  +2.57%  "

$ws.Range("D31").Value = "1.081"
$ws.Range("E31").Value = "  +4.55%  "

$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").Value = "6.010"
$ws.Range("E33").Value = "  +4.70%  "

$ws.Range("D34").Value = "3.615"
$ws.Range("E34").Value = "  -0.13%  "

$ws.Range("D35").Value = "0.02490"
$ws.Range("E35").Value = "  +1.99%  "

$ws.Range("D36").Value = "0.06579"
$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("E37").Value = "  +2.36%  "

$ws.Range("D38").Value = "1.228"
$ws.Range("E38").Value = "  +4.31%  "

$ws.Range("D39").Value = "5.174"
$ws.Range("E39").Value = "  +3.20%  "

$ws.Range("D40").Value = "8.866"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").Value = "0.6533"
$ws.Range("E41").Value = "  +2.90%  "

$ws.Range("D42").Value = "11.64"
$ws.Range("E42").Value = "  +5.15%  "

$ws.Range("D43").Value = "1.243"
$ws.Range("E43").Value = "  +1.07%  "

$ws.Range("D44").Value = "0.6146"
$ws.Range("E44").Value = "  +2.63%  "

$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  +1.72%  "

$ws.Range("D46").Value = "3.746"
$ws.Range("E46").Value = "  +1.75%  "

$ws.Range("D47").Value = "2.080"
$ws.Range("E47").Value = "  +4.40%  "

$ws.Range("E48").Value = "  +3.01%  "

$ws.Range("D49").Value = "124.39"
$ws.Range("E49").Value = "  +1.86%  "

$ws.Range("E50").Value = "  +2.94%  "

$ws.Range("D51").Value = "79.52"
$ws.Range("E51").Value = "  +4.40%  "

$dRange.Style = "Normal"
